$d = $word.ActiveDocument

# Locate the target paragraph that contains the sentence we need to rewrite.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Copy the provided link*") {
        $target = $p.Range
    }
}

$xmlFrag = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="535BE096" w14:textId="5D813B35" w:rsidR="00B90703" w:rsidRDefault="00B90703" w:rsidP="00B90703" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rStyle w:val="edit-btn"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Click on the PIN and copy the provided link that leads to the survey-taker (“user”) login page, either with or without the PIN (if the link without the PIN is chosen, be sure to include the PIN in the email), and email it / them to any entity that you wish to have </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>take</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> that survey. PINs can be used for any number of groups, and by any number of people in a group, but survey results can only be divided by PIN (using group names).</w:t></w:r></w:p>
'@

# Replace the paragraph's content with the restructured runs / proofErr markers.
$target.InsertXML($xmlFrag)

# Re-resolve the (now rebuilt) paragraph and re-apply the "edit-btn" character
# style to its text run(s) -- InsertXML does not round-trip a run-level
# rStyle, so it is re-applied here via the object model (the paragraph-mark's
# rStyle, inside pPr/rPr, round-trips fine and was kept in the XML above).
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Click on the PIN*") {
        $target2 = $p.Range
    }
}
$body = $d.Range($target2.Start, $target2.End - 1)
$body.Style = "edit-btn"

# Move the "_GoBack" bookmark to the start of this paragraph (last edited
# location) -- adding a bookmark named "_GoBack" automatically relocates it
# from its previous position instead of creating a duplicate.
$startPoint = $d.Range($target2.Start, $target2.Start)
$d.Bookmarks.Add("_GoBack", $startPoint)

Write-Host "done"
